# ---------------------------------------------------------------------------
# Applies the two changes captured by the commit diff:
#
#   1. The table on slide 16 gets a new table style
#      ({259F8481-...} -> {75FCA493-588E-4AA3-8AB3-AD83A76391AC}).
#
#   2. The presentation's theme (ppt/theme/theme1.xml, shared by every
#      slide through the slide master) has its 12 theme colors replaced
#      with the stock Office Theme palette (it previously held the
#      "Integral" palette). The theme's font scheme / format scheme are
#      untouched because they were already identical between the two
#      named themes in this deck, so only a:clrScheme differs.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 16 -------------------------------------------

$tableSlide = $p.Slides.Item(16)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shp = $tableSlide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{75FCA493-588E-4AA3-8AB3-AD83A76391AC}")
    }
}

# --- 2. Swap the theme colors for the stock "Office Theme" palette --------

function HexToOle([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Order matches ThemeColorScheme.Item(1..12):
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeThemeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

# Any slide can be used to reach the shared theme part - it is the same
# ppt/theme/theme1.xml backing every slide through the slide master.
$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = HexToOle($officeThemeColors[$i - 1])
}
